# Actualizacion automatica del mapa (2025-10-16 12:53:45)
#
# The underlying data feed dropped the oldest still-pending row (Caso 6036,
# MEDRANO 1715) which shifts every following record up by one row, and a
# brand new record (Caso 2470, DIAZ, CNEL. AV. 2599) was appended at the
# bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the obsolete first data row (row 8 = Caso 6036). Everything
#    below shifts up by one row automatically.
$ws.Rows(8).Delete()

# The row that lands on 8 after the shift (former row 9, Caso 6071) has an
# empty "Observaciones" cell in the source data - make sure it stays a truly
# blank cell rather than an empty string after the shift.
$ws.Range("H8").ClearContents()

# 2) Append the new record as the new last row of the table.
$newRow = 74

$textCols = @("A","B","C","D","E","F","G","H","J","K","L","O","P","Q","R")
foreach ($col in $textCols) {
    $ws.Range($col + $newRow).NumberFormat = "@"
}

$ws.Range("A$newRow").Value = "2470"
$ws.Range("B$newRow").Value = "10/15/2025"
$ws.Range("C$newRow").Value = "DIAZ, CNEL. AV. 2599"
$ws.Range("D$newRow").Value = "14"
$ws.Range("E$newRow").Value = "810371048"
$ws.Range("F$newRow").Value = "NEW"
$ws.Range("G$newRow").Value = "Pendiente"
$ws.Range("H$newRow").Value = "Cambiar"
$ws.Range("I$newRow").Value = 1
$ws.Range("J$newRow").Value = "Cambio"
$ws.Range("K$newRow").Value = "Sin equipos"
$ws.Range("L$newRow").Value = "Pasante"
$ws.Range("M$newRow").Value = -58.405559
$ws.Range("N$newRow").Value = -34.582478
$ws.Range("O$newRow").Value = "Recoleta"
$ws.Range("P$newRow").Value = "Capital Sur"
$ws.Range("Q$newRow").Value = "AGU-N"
$ws.Range("R$newRow").Value = "Fuera de Poligono OVL"
